# Change the column header names in PhotometricParam.xlsx from the old
# LaTeX-flavoured labels (e.g. "$M_{U}$", "$BC_{K}$") to plain-text labels
# (e.g. "Mu", "BCk"), as described in the commit "Change collumn names in
# photometricParam.xlsx".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B1").Value = "Mu"
$ws.Range("C1").Value = "Mb"
$ws.Range("D1").Value = "Mv"
$ws.Range("E1").Value = "Mj"
$ws.Range("F1").Value = "Mh"
$ws.Range("G1").Value = "Mk"
$ws.Range("H1").Value = "(U-B)0"
$ws.Range("I1").Value = "(B-V)0"
$ws.Range("J1").Value = "(J_H)0"
$ws.Range("K1").Value = "(H-K)0"
$ws.Range("L1").Value = "Bcu"
$ws.Range("M1").Value = "BCb"
$ws.Range("N1").Value = "BCv"
$ws.Range("O1").Value = "BCj"
$ws.Range("P1").Value = "BCh"
$ws.Range("Q1").Value = "BCk"

# The shorter plain-text labels no longer need the extra row height the
# old LaTeX labels used, so re-fit the header/body rows (this also drops
# the per-row custom height back to the sheet default, matching the
# saved file).
$ws.Range("A1:A37").EntireRow.AutoFit() | Out-Null

# A few of the new headers are wider than their column, so widen those
# columns to fit the new text (as Excel does automatically).
$ws.Columns.Item(8).ColumnWidth = 13
$ws.Columns.Item(9).ColumnWidth = 12.83
$ws.Columns.Item(10).ColumnWidth = 9
$ws.Columns.Item(11).ColumnWidth = 9.33

# Leave the selection where the edit finished, on the last header cell.
$ws.Range("Q1").Select()
